# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (H) and "Correspond Handback
# DateTime" (K) timestamps for the first data row (the
# 86196a87-d98f-4165-9f6c-55b3f6325547.md entry) on both the "zh-cn" and
# "de-de" report sheets, as produced by a fresh handback-status report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 14:47:53"
$wsZhCn.Range("K2").Value = "2016-08-28 14:48:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 14:47:59"
$wsDeDe.Range("K2").Value = "2016-08-28 14:48:20"
